# project plan update, corresponding budget updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet   # "DRAFT" budget sheet

# UG Research Assistants-Summer (row 12): the increment that used to land in
# Year 5 (I12 = H12*1.02) now lands entirely in Year 4 instead -- Year 4's
# share goes from 33% to 67% of the prior-year figure, and Year 5 drops to 0.
$ws.Range("H12").Formula = "=G12*1.02*0.67"
$ws.Range("I12").Formula = "=H12*1.02*0"

# Drone parts (row 47): Year 4 now gets funded the same as Year 5 (was 0).
$ws.Range("H47").Value = 1000

# Leave the selection where the editor ended up working.
$ws.Range("D13").Select()
